# Add a new "source_uri" column to the "Sample" sheet, inserted right
# before the existing "id" column (which currently lives in column F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# Insert a new column at F, shifting id/name/description (F,G,H) to G,H,I
$ws.Range("F1").EntireColumn.Insert()

$ws.Range("F1").Value = "source_uri"
